# "Generate Report for Archive"
# The localization-status report is regenerated: the e2e test file
# 081d8076-ff0e-4c03-98b7-51c2af6d19b4.md has moved from "Ready for
# handoff" to "In Translation", and the report rows for the three
# "In Translation" files are re-sorted alphabetically by file name
# (081d8076 < 12c51228 < 440d60bd). This shifts the data that used to
# live in rows 3-5 of each sheet (Overview, zh-cn, de-de) around; the
# rest of the rows (2, 6, 7) are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview": columns A-G, rows 3-5 get new content.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = "081d8076-ff0e-4c03-98b7-51c2af6d19b4.md"
$ws1.Range("B3").Value = "e2e\081d8076-ff0e-4c03-98b7-51c2af6d19b4.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"
$ws1.Range("G3").Value = "2016-08-21 14:49:15"

$ws1.Range("A4").Value = "12c51228-8784-4454-bc7c-ae7d0be05400.md"
$ws1.Range("B4").Value = "e2e\12c51228-8784-4454-bc7c-ae7d0be05400.md"
$ws1.Range("C4").Value = ".md"
$ws1.Range("E4").Value = "In Translation"
$ws1.Range("F4").Value = "In Translation"
$ws1.Range("G4").Value = "2016-08-21 14:48:39"

$ws1.Range("A5").Value = "440d60bd-8fb4-4c83-9f59-c76cf27d4766.md"
$ws1.Range("B5").Value = "e2e\440d60bd-8fb4-4c83-9f59-c76cf27d4766.md"
$ws1.Range("C5").Value = ".md"
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("G5").Value = "2016-08-21 14:48:39"

# Rebuild the B2:B7 hyperlinks so the displayed text matches the new
# row order. The underlying relationship targets (Address) are kept
# exactly as they were wired before the edit (rId3->12c51228's URL,
# rId4->440d60bd's URL, rId5->081d8076's URL, ...) - only which row
# (and therefore which display text) uses each target changes.
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/212f18de0bcef964ae26ef88818e28b801442ff9/e2e/26ca0bb9-56a4-45f3-b058-f36de3211cfd.md", "", "", "e2e\26ca0bb9-56a4-45f3-b058-f36de3211cfd.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1083303a25c6742b69c4812dfc3ad7dac56707c4/e2e/12c51228-8784-4454-bc7c-ae7d0be05400.md", "", "", "e2e\081d8076-ff0e-4c03-98b7-51c2af6d19b4.md")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1083303a25c6742b69c4812dfc3ad7dac56707c4/e2e/440d60bd-8fb4-4c83-9f59-c76cf27d4766.md", "", "", "e2e\12c51228-8784-4454-bc7c-ae7d0be05400.md")
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeff86f6e3378e721e770562e239943d670f473c/e2e/081d8076-ff0e-4c03-98b7-51c2af6d19b4.md", "", "", "e2e\440d60bd-8fb4-4c83-9f59-c76cf27d4766.md")
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0dbfaac9784a0b5c3876d6304e74fb79bc289b6/e2e/796762e4-7f92-41be-85a6-414a4b2e9726.md", "", "", "e2e\796762e4-7f92-41be-85a6-414a4b2e9726.md")
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e02e7f33666a8254ffc39b2b2a1df6386e4ed92/e2e/ac0b8092-787e-4c1c-9179-d816e2c0177f.md", "", "", "e2e\ac0b8092-787e-4c1c-9179-d816e2c0177f.md")

# ---------------------------------------------------------------
# Sheet "zh-cn": columns A-P, rows 3-5 get new content.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = "081d8076-ff0e-4c03-98b7-51c2af6d19b4.md"
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("G3").Value = "081d8076-ff0e-4c03-98b7-51c2af6d19b4.04e13ae4ae309559fc6514a43188a78e4228f249.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-21 14:49:11"

$ws2.Range("A4").Value = "12c51228-8784-4454-bc7c-ae7d0be05400.md"
$ws2.Range("C4").Value = "In Translation"
$ws2.Range("G4").Value = "12c51228-8784-4454-bc7c-ae7d0be05400.b7b4aad2f0dede85f7d0241396fa6275f6695723.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-21 14:48:34"

$ws2.Range("A5").Value = "440d60bd-8fb4-4c83-9f59-c76cf27d4766.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("G5").Value = "440d60bd-8fb4-4c83-9f59-c76cf27d4766.76c25fd71230fc8836122e6aa64846927a4207d9.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-21 14:48:34"

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/212f18de0bcef964ae26ef88818e28b801442ff9/e2e/26ca0bb9-56a4-45f3-b058-f36de3211cfd.md", "", "", "26ca0bb9-56a4-45f3-b058-f36de3211cfd.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d49e0b9cbac704788053ff6e0ef2a557fa50b380/e2e/26ca0bb9-56a4-45f3-b058-f36de3211cfd.md", "", "", "26ca0bb9-56a4-45f3-b058-f36de3211cfd.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1083303a25c6742b69c4812dfc3ad7dac56707c4/e2e/12c51228-8784-4454-bc7c-ae7d0be05400.md", "", "", "081d8076-ff0e-4c03-98b7-51c2af6d19b4.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1083303a25c6742b69c4812dfc3ad7dac56707c4/e2e/440d60bd-8fb4-4c83-9f59-c76cf27d4766.md", "", "", "12c51228-8784-4454-bc7c-ae7d0be05400.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeff86f6e3378e721e770562e239943d670f473c/e2e/081d8076-ff0e-4c03-98b7-51c2af6d19b4.md", "", "", "440d60bd-8fb4-4c83-9f59-c76cf27d4766.md")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0dbfaac9784a0b5c3876d6304e74fb79bc289b6/e2e/796762e4-7f92-41be-85a6-414a4b2e9726.md", "", "", "796762e4-7f92-41be-85a6-414a4b2e9726.md")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e02e7f33666a8254ffc39b2b2a1df6386e4ed92/e2e/ac0b8092-787e-4c1c-9179-d816e2c0177f.md", "", "", "ac0b8092-787e-4c1c-9179-d816e2c0177f.md")

# ---------------------------------------------------------------
# Sheet "de-de": columns A-P, rows 3-5 get new content.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = "081d8076-ff0e-4c03-98b7-51c2af6d19b4.md"
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("G3").Value = "081d8076-ff0e-4c03-98b7-51c2af6d19b4.04e13ae4ae309559fc6514a43188a78e4228f249.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-21 14:49:15"

$ws3.Range("A4").Value = "12c51228-8784-4454-bc7c-ae7d0be05400.md"
$ws3.Range("C4").Value = "In Translation"
$ws3.Range("G4").Value = "12c51228-8784-4454-bc7c-ae7d0be05400.b7b4aad2f0dede85f7d0241396fa6275f6695723.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-21 14:48:39"

$ws3.Range("A5").Value = "440d60bd-8fb4-4c83-9f59-c76cf27d4766.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("G5").Value = "440d60bd-8fb4-4c83-9f59-c76cf27d4766.76c25fd71230fc8836122e6aa64846927a4207d9.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-21 14:48:39"

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/212f18de0bcef964ae26ef88818e28b801442ff9/e2e/26ca0bb9-56a4-45f3-b058-f36de3211cfd.md", "", "", "26ca0bb9-56a4-45f3-b058-f36de3211cfd.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bbd49b33d32c09e950a9446500e5e86dbb61ee21/e2e/26ca0bb9-56a4-45f3-b058-f36de3211cfd.md", "", "", "26ca0bb9-56a4-45f3-b058-f36de3211cfd.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1083303a25c6742b69c4812dfc3ad7dac56707c4/e2e/12c51228-8784-4454-bc7c-ae7d0be05400.md", "", "", "081d8076-ff0e-4c03-98b7-51c2af6d19b4.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1083303a25c6742b69c4812dfc3ad7dac56707c4/e2e/440d60bd-8fb4-4c83-9f59-c76cf27d4766.md", "", "", "12c51228-8784-4454-bc7c-ae7d0be05400.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeff86f6e3378e721e770562e239943d670f473c/e2e/081d8076-ff0e-4c03-98b7-51c2af6d19b4.md", "", "", "440d60bd-8fb4-4c83-9f59-c76cf27d4766.md")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0dbfaac9784a0b5c3876d6304e74fb79bc289b6/e2e/796762e4-7f92-41be-85a6-414a4b2e9726.md", "", "", "796762e4-7f92-41be-85a6-414a4b2e9726.md")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e02e7f33666a8254ffc39b2b2a1df6386e4ed92/e2e/ac0b8092-787e-4c1c-9179-d816e2c0177f.md", "", "", "ac0b8092-787e-4c1c-9179-d816e2c0177f.md")
